# cap nhat tai lieu: lich trinh, dac ta yeu cau, cap nhat ke hoach
$d = $word.ActiveDocument

# --- 1. "6. Nguon luc" table: rename existing "Khoa" entries to "Khanh" ---
$d.Content.Find.Execute("Khoa", $true, $true, $false, $false, $false, $true, 1, $false, "Khanh", 2) | Out-Null

# --- 2. Schedule table: design row gains a third member ---
$d.Content.Find.Execute("Ngôn, Kiệt", $true, $false, $false, $false, $false, $true, 1, $false, "Ngôn, Kiệt, Khoa", 2) | Out-Null

# --- 3. Schedule table: "Lap trinh" task now also covers interface coding ---
$d.Content.Find.Execute("Lập trình (Code các chức năng chính)", $true, $false, $false, $false, $false, $true, 1, $false, "Lập trình (Code giao diện các chức năng chính)", 2) | Out-Null

# --- 4. Schedule table: programming task owners expanded ---
$d.Content.Find.Execute("GHuy, Hộp", $true, $false, $false, $false, $false, $true, 1, $false, "Vinh Huy, Nghĩa, Gia Huy, Hộp", 2) | Out-Null

# --- 5. "6. Nguon luc" table: insert a new row for Khoa (UI/UX Designer) right after Kiet's row ---
$t = $d.Tables(2)
$kietRow = $t.Rows(4)
$newRow = $t.Rows.Add($t.Rows(5))
$newRow.Cells(1).Range.Text = "Khoa"
$newRow.Cells(2).Range.Text = "UI/UX Designer"
$newRow.Cells(3).Range.Text = "Thiết kế giao diện người dùng, trải nghiệm người dùng."

# --- 6. Schedule table column widths updated ---
$sched = $d.Tables(1)
$sched.Columns(1).Width = 175.35
$sched.Columns(2).Width = 114
$sched.Columns(3).Width = 71.1
$sched.Columns(4).Width = 70.75
